$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 2.469524
$ws.Range("H2").Value = 7.408571999999999
$ws.Range("I2").Value = 0.006775482240913427
$ws.Range("J2").Value = 0.006775482240913427
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.731629
$ws.Range("N2").Value = 8.194887
$ws.Range("O2").Value = 0.5547800938501829
$ws.Range("P2").Value = 0.554780093850183
$ws.Range("Q2").Value = 6.745823374595999
$ws.Range("R2").Value = 60.71241037136399
$ws.Range("S2").Value = 0.003758902673494199
$ws.Range("T2").Value = 0.003758902673494199

# Row 3
$ws.Range("G3").Value = 2.469524
$ws.Range("H3").Value = 7.408571999999999
$ws.Range("I3").Value = 0.006775482240913427
$ws.Range("J3").Value = 0.006775482240913427
$ws.Range("M3").Value = 0.06813733333333333
$ws.Range("O3").Value = 0.01383834927121065
$ws.Range("P3").Value = 0.01383834927121065
$ws.Range("Q3").Value = 0.1682667799626666
$ws.Range("R3").Value = 1.514401019664
$ws.Range("S3").Value = 0.00009376148973064499
$ws.Range("T3").Value = 0.00009376148973064502

# Row 4
$ws.Range("G4").Value = 2.469524
$ws.Range("H4").Value = 7.408571999999999
$ws.Range("I4").Value = 0.006775482240913427
$ws.Range("J4").Value = 0.006775482240913427
$ws.Range("M4").Value = 2.124038666666666
$ws.Range("N4").Value = 6.372115999999999
$ws.Range("O4").Value = 0.4313815568786064
$ws.Range("P4").Value = 0.4313815568786064
$ws.Range("Q4").Value = 5.245364464261332
$ws.Range("R4").Value = 47.20828017835199
$ws.Range("S4").Value = 0.002922818077688583
$ws.Range("T4").Value = 0.002922818077688583

# Row 5
$ws.Range("H5").Value = 988.862762
$ws.Range("I5").Value = 0.9043607975506752
$ws.Range("J5").Value = 0.9043607975506752
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.731629
$ws.Range("N5").Value = 8.194887
$ws.Range("O5").Value = 0.5547800938501829
$ws.Range("P5").Value = 0.554780093850183
$ws.Range("Q5").Value = 900.4020658997659
$ws.Range("R5").Value = 8103.618593097894
$ws.Range("S5").Value = 0.5017213681395899
$ws.Range("T5").Value = 0.50172136813959

# Row 6
$ws.Range("H6").Value = 988.862762
$ws.Range("I6").Value = 0.9043607975506752
$ws.Range("J6").Value = 0.9043607975506752
$ws.Range("M6").Value = 0.06813733333333333
$ws.Range("O6").Value = 0.01383834927121065
$ws.Range("P6").Value = 0.01383834927121065
$ws.Range("S6").Value = 0.01251486058369686
$ws.Range("T6").Value = 0.01251486058369687

# Row 7
$ws.Range("H7").Value = 988.862762
$ws.Range("I7").Value = 0.9043607975506752
$ws.Range("J7").Value = 0.9043607975506752
$ws.Range("M7").Value = 2.124038666666666
$ws.Range("N7").Value = 6.372115999999999
$ws.Range("O7").Value = 0.4313815568786064
$ws.Range("P7").Value = 0.4313815568786064
$ws.Range("Q7").Value = 700.1275808382658
$ws.Range("R7").Value = 6301.148227544391
$ws.Range("S7").Value = 0.3901245688273884
$ws.Range("T7").Value = 0.3901245688273884

# Row 8
$ws.Range("G8").Value = 32.38899933333333
$ws.Range("H8").Value = 97.16699799999999
$ws.Range("I8").Value = 0.08886372020841134
$ws.Range("J8").Value = 0.08886372020841135
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.731629
$ws.Range("N8").Value = 8.194887
$ws.Range("O8").Value = 0.5547800938501829
$ws.Range("P8").Value = 0.554780093850183
$ws.Range("Q8").Value = 88.47472985991399
$ws.Range("R8").Value = 796.2725687392259
$ws.Range("S8").Value = 0.04929982303709884
$ws.Range("T8").Value = 0.04929982303709885

# Row 9
$ws.Range("G9").Value = 32.38899933333333
$ws.Range("H9").Value = 97.16699799999999
$ws.Range("I9").Value = 0.08886372020841134
$ws.Range("J9").Value = 0.08886372020841135
$ws.Range("M9").Value = 0.06813733333333333
$ws.Range("O9").Value = 0.01383834927121065
$ws.Range("P9").Value = 0.01383834927121065
$ws.Range("Q9").Value = 2.206900043908444
$ws.Range("R9").Value = 19.862100395176
$ws.Range("S9").Value = 0.001229727197783136
$ws.Range("T9").Value = 0.001229727197783136

# Row 10
$ws.Range("G10").Value = 32.38899933333333
$ws.Range("H10").Value = 97.16699799999999
$ws.Range("I10").Value = 0.08886372020841134
$ws.Range("J10").Value = 0.08886372020841135
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.124038666666666
$ws.Range("N10").Value = 6.372115999999999
$ws.Range("O10").Value = 0.4313815568786064
$ws.Range("P10").Value = 0.4313815568786064
$ws.Range("Q10").Value = 68.79548695864088
$ws.Range("R10").Value = 619.1593826277679
$ws.Range("S10").Value = 0.03833416997352936
$ws.Range("T10").Value = 0.03833416997352937
